# Commit: "atualiza paper in press"
# Updates hrefs (double-backslash -> single backslash, adds scheme/www),
# turns the botanicaamazonica.wiki.br link in F6 into a real hyperlink whose
# visible "http://www.botanicaamazonica" substring is colored blue, fixes a
# couple of column widths, and re-points the active sheet/tab/selection from
# "grants" back to "experience".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experience")

# --- plain text fixes (single backslash \href, scheme/www added where needed) ---

$ws.Range("F2").Value = 'Evento ocorrido durante a \href{https://www.sncticet.ufam.edu.br}{XV Semana Nacional de Ciência e Tecnologia de Itacoatiara (SNCT-ITA)};Duração de quatro horas, divididas em dois dias, duas horas cada; Este minicurso teve a intenção de abordar simplificadamente a produção de mapas de distribuição de espécies utilizando a linguagem R'

$ws.Range("F7").Value = 'Professores: Dr. Reinaldo Imbrozio Barbosa (INPA), Lidiany Carvalho (UFRR);Auxiliei discentes em lidar com o ambiente R;Contribuí ativamente para o ensino do curso através de reuniões com o professor;Criei um \href{http://www.botanicaamazonica.wiki.br/labotam/doku.php?id=alunos:r.perdiz:disciplina:inicio}{sítio web} para auxiliar os discentes no aprendizado do R'

$ws.Range("F9").Value = 'Event that took place during the \href{https://www.sncticet.ufam.edu.br}{XV Semana Nacional de Ciência e Tecnologia de Itacoatiara (SNCT-ITA)};Duration of four hours, divided into two days, two hours each;This short course had the intention of teaching in a simplified way how to produce species distribution maps using the R language'

$ws.Range("F13").Value = 'Teacher: Dr. Alberto Vicentini (INPA);Check the \href{http://www.botanicaamazonica.wiki.br/labotam/doku.php?id=disciplinas:bot89:inicio}{website}'

$ws.Range("F14").Value = 'Teacher: Dr. Reinaldo Imbrozio Barbosa (INPA), Lidiany Carvalho (UFRR);Assisted students in dealing with R environment;Actively contributed to lecture course design through meetings with instructor;Created a \href{http://www.botanicaamazonica.wiki.br/labotam/doku.php?id=alunos:r.perdiz:disciplina:inicio}{website} to help students in learning R'

# --- F6: rich text with a blue-colored run + a real hyperlink ---

# Apply the blue Arial/10 font to a scratch cell first so that font gets
# registered in the workbook's shared font table, then clear the scratch
# cell again so it leaves no trace in the sheet.
$ws.Range("Z100").Value = "x"
$ws.Range("Z100").Font.Color = 16711680
$ws.Range("Z100").Font.Size = 10
$ws.Range("Z100").Font.Name = "Arial"
$ws.Range("Z100").Clear()

$f6text = 'Professor: Dr. Alberto Vicentini (INPA);Atuei como monitor por dois anos consecutivos;Criei um \href{http://www.botanicaamazonica.wiki.br/labotam/doku.php?id=disciplinas:bot89:inicio}{sítio web} para auxiliar os discentes no aprendizado da disciplina'
$ws.Range("F6").Value = $f6text

# run 1: "Professor: ... \href{"  -- length 101, default color
$ws.Range("F6").Characters(1, 101).Font.ColorIndex = -4105
# run 2: "http://www.botanicaamazonica" -- length 28, blue
$ws.Range("F6").Characters(102, 28).Font.Color = 16711680
# run 3: ".wiki.br/labotam/...{sítio web} ..." -- remainder, default color
$ws.Range("F6").Characters(130, 121).Font.ColorIndex = -4105

$ws.Hyperlinks.Add($ws.Range("F6"), "http://www.botanicaamazonica.wiki.br/labotam/doku.php?id=disciplinas:bot89:inicio", "", "", "http://www.botanicaamazonica")

# --- column width tweaks ---

$ws.Columns.Item(6).ColumnWidth = 62

$degree = $wb.Worksheets.Item("degree")
$degree.Columns.Item(6).ColumnWidth = 52.5

# --- active sheet / tab / selection: back to "experience" (was "grants") ---

$ws.Activate()
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 6
